$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# --- Row 23/24 swap: Toncoin moves to row 23, BinanceUSD moves to row 24 ---
$ws.Range("B23").Value = "Toncoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D23") "1.937"
$ws.Range("E23").Value = "  +4.47%  "

$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws.Range("D24") "1.003"
$ws.Range("E24").Value = "  -0.55%  "

# --- Price / Volume updates ---
$ws.Range("D2").Value = "25.973.09"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "1.638.52"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("E4").Value = "  -0.49%  "
Set-TextValue $ws.Range("D5") "215.03"
$ws.Range("E5").Value = "  -0.17%  "
Set-TextValue $ws.Range("D6") "0.5141"
$ws.Range("E6").Value = "  +0.69%  "
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("E8").Value = "  -0.02%  "
Set-TextValue $ws.Range("D9") "0.06350"
$ws.Range("E9").Value = "  -1.06%  "
Set-TextValue $ws.Range("D10") "19.77"
$ws.Range("E10").Value = "  +0.29%  "
Set-TextValue $ws.Range("D11") "0.07765"
$ws.Range("E11").Value = "  -0.14%  "
Set-TextValue $ws.Range("D12") "4.273"
$ws.Range("D13").Value = "1.621.78"
$ws.Range("E13").Value = "  -1.79%  "
$ws.Range("E14").Value = "  -0.39%  "
$ws.Range("D15").Value = "0.0₅7748"
$ws.Range("E15").Value = "  -1.89%  "
Set-TextValue $ws.Range("D16") "64.33"
$ws.Range("E16").Value = "  -0.97%  "
$ws.Range("D17").Value = "25.988.27"
$ws.Range("E17").Value = "  +0.00%  "
Set-TextValue $ws.Range("D18") "0.9982"
$ws.Range("E18").Value = "  -0.79%  "
Set-TextValue $ws.Range("D19") "197.24"
$ws.Range("E19").Value = "  -0.23%  "
Set-TextValue $ws.Range("D20") "4.433"
$ws.Range("E20").Value = "  -0.01%  "
Set-TextValue $ws.Range("D21") "9.917"
$ws.Range("E21").Value = "  -1.19%  "
Set-TextValue $ws.Range("D22") "6.081"
$ws.Range("E22").Value = "  +0.32%  "
Set-TextValue $ws.Range("D25") "142.09"
$ws.Range("E25").Value = "  +1.11%  "
Set-TextValue $ws.Range("D26") "0.1234"
$ws.Range("E26").Value = "  +7.53%  "
Set-TextValue $ws.Range("D27") "6.833"
$ws.Range("E27").Value = "  -0.89%  "
Set-TextValue $ws.Range("D28") "15.57"
$ws.Range("E28").Value = "  -1.21%  "
Set-TextValue $ws.Range("D29") "1.239"
$ws.Range("E29").Value = "  -0.04%  "
Set-TextValue $ws.Range("D30") "0.04842"
Set-TextValue $ws.Range("D31") "3.281"
$ws.Range("E31").Value = "  +0.20%  "
Set-TextValue $ws.Range("D32") "3.210"
$ws.Range("E32").Value = "  +0.25%  "
Set-TextValue $ws.Range("D33") "1.532"
$ws.Range("E33").Value = "  -0.75%  "
$ws.Range("E34").Value = "  +0.33%  "
Set-TextValue $ws.Range("D35") "0.9116"
$ws.Range("E35").Value = "  +1.88%  "
Set-TextValue $ws.Range("D36") "2.566"
$ws.Range("E36").Value = "  -0.78%  "
Set-TextValue $ws.Range("D37") "0.5538"
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("D38").Value = "1.105.24"
$ws.Range("E38").Value = "  -2.43%  "
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("E40").Value = "  -0.51%  "
Set-TextValue $ws.Range("D41") "2.522"
$ws.Range("E41").Value = "  -1.61%  "
Set-TextValue $ws.Range("D42") "5.556"
$ws.Range("E42").Value = "  -1.83%  "
Set-TextValue $ws.Range("D43") "0.8053"
$ws.Range("E43").Value = "  -1.14%  "
Set-TextValue $ws.Range("D44") "99.19"
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("E45").Value = "  -2.40%  "
$ws.Range("D46").Value = "1.778.17"
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("E47").Value = "  +0.05%  "
Set-TextValue $ws.Range("D48") "55.02"
$ws.Range("E48").Value = "  -0.57%  "
Set-TextValue $ws.Range("D49") "0.9979"
$ws.Range("E49").Value = "  -0.76%  "
Set-TextValue $ws.Range("D50") "0.05210"
$ws.Range("E50").Value = "  +2.32%  "
Set-TextValue $ws.Range("D51") "7.493"
$ws.Range("E51").Value = "  +1.27%  "
